# Add Minute5/Second5/Rep5 data (columns S, T, U) for rows 2-19
# on the "Score" worksheet, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

$data = @{
    2  = @(4, 0, 2)
    3  = @(2, 56, 2)
    4  = @(6, 13, 2)
    5  = @(5, 7, 2)
    6  = @(4, 29, 2)
    7  = @(3, 47, 2)
    8  = @(4, 32, 2)
    9  = @(3, 26, 2)
    10 = @(4, 33, 2)
    11 = @(3, 59, 2)
    12 = @(4, 13, 2)
    13 = @(3, 39, 2)
    14 = @(2, 59, 2)
    15 = @(4, 38, 2)
    16 = @(4, 16, 2)
    17 = @(4, 24, 2)
    18 = @(4, 34, 2)
    19 = @(4, 25, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 19).Value = $values[0]  # Column S
    $ws.Cells.Item($row, 20).Value = $values[1]  # Column T
    $ws.Cells.Item($row, 21).Value = $values[2]  # Column U
}

# Update view: scroll so column B is the top-left, and select V4
$ws.Activate()
$ws.Range("V4").Select()
$excel.ActiveWindow.ScrollColumn = 2
